$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# 1) Paragraph 1 ("Supplement S5 File"):
#    - remove the existing _GoBack bookmark
#    - change the "5" run's text to "8"
# -------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$p1 = $d.Paragraphs.Item(1).Range
$p1Start = $p1.Start
$prefix1 = "Supplement S"
$digitIdx1 = $p1Start + $prefix1.Length

# Bracket the single-character "5" run with temporary bookmarks so that
# editing its text does not cause it to merge with the neighbouring runs
# (which share identical run formatting).
$d.Bookmarks.Add("ZZTMP1A", $d.Range($digitIdx1, $digitIdx1)) | Out-Null
$d.Bookmarks.Add("ZZTMP1B", $d.Range($digitIdx1 + 1, $digitIdx1 + 1)) | Out-Null
$d.Range($digitIdx1, $digitIdx1 + 1).Text = "8"
$d.Bookmarks("ZZTMP1A").Delete()
$d.Bookmarks("ZZTMP1B").Delete()

# -------------------------------------------------------------------------
# 2) Paragraph 2 ("S5 File. Hardware details."):
#    - change the "5" run's text to "8"
#    - re-insert the _GoBack bookmark right after that run (between the
#      "8" run and the " File." run)
# -------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2).Range
$p2Start = $p2.Start
$prefix2 = "S"
$digitIdx2 = $p2Start + $prefix2.Length

$d.Bookmarks.Add("ZZTMP2A", $d.Range($digitIdx2, $digitIdx2)) | Out-Null
$d.Bookmarks.Add("ZZTMP2B", $d.Range($digitIdx2 + 1, $digitIdx2 + 1)) | Out-Null
$d.Range($digitIdx2, $digitIdx2 + 1).Text = "8"
$d.Bookmarks("ZZTMP2A").Delete()
$d.Bookmarks("ZZTMP2B").Delete()

# The bookmark boundary right after the (now) "8" run is still at $digitIdx2 + 1
$d.Bookmarks.Add("_GoBack", $d.Range($digitIdx2 + 1, $digitIdx2 + 1)) | Out-Null

# -------------------------------------------------------------------------
# 3) Paragraph 3 ("All image processing was completed on the SeaWulf ..."):
#    - split the opening run into three runs:
#        "All image processing"
#        " and model training and validation"
#        " was completed on the "
# -------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3).Range
$p3Start = $p3.Start
$oldOpening = "All image processing was completed on the "
$part1 = "All image processing"
$part2 = " and model training and validation"
$part3 = " was completed on the "

$openRange = $d.Range($p3Start, $p3Start + $oldOpening.Length)
$openRange.Text = $part1 + $part2 + $part3

$splitA = $p3Start + $part1.Length
$splitB = $splitA + $part2.Length

$d.Bookmarks.Add("ZZSPLITA", $d.Range($splitA, $splitA)) | Out-Null
$d.Bookmarks.Add("ZZSPLITB", $d.Range($splitB, $splitB)) | Out-Null

# The first run ("All image processing") inherited a stray
# xml:space="preserve" flag from the original (space-terminated) run even
# though its own text has no leading/trailing space. Force the engine to
# recompute that flag by dirtying the run's text and setting it back.
$r1 = $d.Range($p3Start, $splitA)
$r1.Text = $r1.Text + "Z"
$r1b = $d.Range($p3Start, $splitA + 1)
$r1b.Text = $part1

$d.Bookmarks("ZZSPLITA").Delete()
$d.Bookmarks("ZZSPLITB").Delete()
